$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 34.666668
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H113").Value = 4638.3687
$ws.Range("I113").Value = 4870.3335
$ws.Range("J113").Value = 4240.7144
$ws.Range("K113").Value = 4870.3335
$ws.Range("L113").Value = 4240.7144
$ws.Range("M113").Value = -1616.3335
$ws.Range("N113").Value = -10748.7144
$ws.Range("H132").Value = 109922.28
$ws.Range("I132").Value = 265108.88
$ws.Range("K132").Value = 795326.64
$ws.Range("M132").Value = -792796.64
$ws.Range("H137").Value = 2419.84
$ws.Range("I137").Value = 1708.2354
$ws.Range("K137").Value = 5124.706200000001
$ws.Range("M137").Value = -2574.706200000001
$ws.Range("H138").Value = 5477.66
$ws.Range("I138").Value = 1145.8
$ws.Range("J138").Value = 6485.07
$ws.Range("K138").Value = 3437.4
$ws.Range("L138").Value = 19455.21
$ws.Range("M138").Value = 1702.6
$ws.Range("N138").Value = -29735.21

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3712.1384
$ws.Range("I32").Value = 2009.2322
$ws.Range("K32").Value = 2009.2322
$ws.Range("M32").Value = -1722.2322
$ws.Range("H61").Value = 23172.21
$ws.Range("I61").Value = 20468
$ws.Range("J61").Value = 26176.889
$ws.Range("K61").Value = 20468
$ws.Range("L61").Value = 26176.889
$ws.Range("M61").Value = -20256
$ws.Range("N61").Value = -26600.889
$ws.Range("H74").Value = 11907495
$ws.Range("I74").Value = 17858456
$ws.Range("K74").Value = 17858456
$ws.Range("M74").Value = -17857582
$ws.Range("H77").Value = 11907495
$ws.Range("I77").Value = 17858456
$ws.Range("K77").Value = 89292280
$ws.Range("M77").Value = -89287912
$ws.Range("H122").Value = 3850.7083
$ws.Range("I122").Value = 2829.3809
$ws.Range("K122").Value = 8488.1427
$ws.Range("M122").Value = -6038.1427
$ws.Range("H132").Value = 18500.62
$ws.Range("I132").Value = 17343.568
$ws.Range("K132").Value = 52030.704
$ws.Range("M132").Value = -49500.704
$ws.Range("H136").Value = 23172.21
$ws.Range("I136").Value = 20468
$ws.Range("J136").Value = 26176.889
$ws.Range("K136").Value = 61404
$ws.Range("L136").Value = 78530.667
$ws.Range("M136").Value = -58854
$ws.Range("N136").Value = -83630.667

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 70000
$ws.Range("J103").Value = 70000
$ws.Range("L103").Value = 70000
$ws.Range("N103").Value = -72344
$ws.Range("H107").Value = 55559670
$ws.Range("I107").Value = 74078410
$ws.Range("J107").Value = 3469
$ws.Range("K107").Value = 74078410
$ws.Range("L107").Value = 3469
$ws.Range("M107").Value = -74076490
$ws.Range("N107").Value = -7309
$ws.Range("H134").Value = 1827.3043
$ws.Range("I134").Value = 1586.381
$ws.Range("K134").Value = 4759.143
$ws.Range("M134").Value = -2224.143

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28574168
$ws.Range("I31").Value = 55557450
$ws.Range("J31").Value = 3634.8823
$ws.Range("K31").Value = 55557450
$ws.Range("L31").Value = 3634.8823
$ws.Range("M31").Value = -55557155
$ws.Range("N31").Value = -4224.8823
$ws.Range("H34").Value = 28574168
$ws.Range("I34").Value = 55557450
$ws.Range("J34").Value = 3634.8823
$ws.Range("K34").Value = 55557450
$ws.Range("L34").Value = 3634.8823
$ws.Range("M34").Value = -55557248
$ws.Range("N34").Value = -4038.8823
$ws.Range("H132").Value = 36043896
$ws.Range("I132").Value = 39218276
$ws.Range("K132").Value = 117654828
$ws.Range("M132").Value = -117652298
$ws.Range("H134").Value = 1774.7693
$ws.Range("I134").Value = 1631.0834
$ws.Range("K134").Value = 4893.2502
$ws.Range("M134").Value = -2358.2502

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 871.1818
$ws.Range("I92").Value = 799.8570999999999
$ws.Range("J92").Value = 996
$ws.Range("K92").Value = 2399.5713
$ws.Range("L92").Value = 2988
$ws.Range("M92").Value = -1151.5713
$ws.Range("N92").Value = -5484
$ws.Range("H94").Value = 7933.375
$ws.Range("J94").Value = 7933.857
$ws.Range("L94").Value = 23801.571
$ws.Range("N94").Value = -25153.571
$ws.Range("H96").Value = 9000
$ws.Range("I96").Value = 5000
$ws.Range("K96").Value = 15000
$ws.Range("M96").Value = -12941
$ws.Range("H97").Value = 298.5
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 298.5
$ws.Range("K97").Value = 0
$ws.Range("L97").ClearContents()
$ws.Range("M97").Value = 895.5
$ws.Range("N97").Value = -1887.5
$ws.Range("H98").Value = 1422.8889
$ws.Range("I98").Value = 990
$ws.Range("K98").Value = 2970
$ws.Range("M98").Value = -1472
$ws.Range("H99").Value = 3455.7144
$ws.Range("I99").Value = 1047.5
$ws.Range("J99").Value = 6666.6665
$ws.Range("K99").Value = 3142.5
$ws.Range("L99").Value = 19999.9995
$ws.Range("M99").Value = -896.5
$ws.Range("N99").Value = -24491.9995
$ws.Range("H100").Value = 5949.5
$ws.Range("J100").Value = 5949.5
$ws.Range("L100").Value = 17848.5
$ws.Range("N100").Value = -19470.5
$ws.Range("H101").Value = 250008260
$ws.Range("J101").Value = 250008260
$ws.Range("L101").Value = 750024780
$ws.Range("N101").Value = -750029648
$ws.Range("H102").Value = 18505
$ws.Range("J102").Value = 18505
$ws.Range("L102").Value = 55515
$ws.Range("N102").Value = -60383
$ws.Range("H103").Value = 960
$ws.Range("I103").Value = 800
$ws.Range("K103").Value = 2400
$ws.Range("M103").Value = -1521
$ws.Range("H104").Value = 9800
$ws.Range("I104").Value = 5000
$ws.Range("K104").Value = 15000
$ws.Range("M104").Value = -12379
$ws.Range("H110").Value = 24999
$ws.Range("I110").Value = 24999
$ws.Range("K110").Value = 74997
$ws.Range("M110").Value = -70907
$ws.Range("H115").Value = 12000
$ws.Range("I115").Value = 4000
$ws.Range("K115").Value = 12000
$ws.Range("M115").Value = -10825
$ws.Range("H132").Value = 2980.5454
$ws.Range("I132").Value = 700
$ws.Range("K132").Value = 6300
$ws.Range("M132").Value = -3770

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 200000
$ws.Range("I49").Value = 200000
$ws.Range("K49").Value = 200000
$ws.Range("M49").Value = -199816
$ws.Range("H70").Value = 6810791.5
$ws.Range("I70").Value = 15878681
$ws.Range("K70").Value = 15878681
$ws.Range("M70").Value = -15878411
$ws.Range("H73").Value = 6810791.5
$ws.Range("I73").Value = 15878681
$ws.Range("K73").Value = 15878681
$ws.Range("M73").Value = -15877745
$ws.Range("H80").Value = 1393512.6
$ws.Range("I80").Value = 2083031.6
$ws.Range("J80").Value = 14474.75
$ws.Range("K80").Value = 2083031.6
$ws.Range("L80").Value = 14474.75
$ws.Range("M80").Value = -2082033.6
$ws.Range("N80").Value = -16470.75
$ws.Range("H83").Value = 1393512.6
$ws.Range("I83").Value = 2083031.6
$ws.Range("J83").Value = 14474.75
$ws.Range("K83").Value = 10415158
$ws.Range("L83").Value = 72373.75
$ws.Range("M83").Value = -10410166
$ws.Range("N83").Value = -82357.75
$ws.Range("H126").Value = 4054.9
$ws.Range("I126").Value = 2073.2666
$ws.Range("K126").Value = 6219.7998
$ws.Range("M126").Value = -3749.7998
$ws.Range("H132").Value = 2983.625
$ws.Range("I132").Value = 2919.675
$ws.Range("K132").Value = 8759.025000000001
$ws.Range("M132").Value = -6229.025000000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3331.0908
$ws.Range("I126").Value = 3265.2
$ws.Range("K126").Value = 9795.599999999999
$ws.Range("M126").Value = -7325.599999999999
